# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Row number -> new F-column value
$updates = @{
    2  = 170
    3  = 427
    4  = 12390
    5  = 1270
    6  = 142
    9  = 159
    10 = 193
    16 = 373
    17 = 3529
    18 = 92
    19 = 939
    22 = 42
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
